$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a stray 3-blank-row gap between row 36 (MessageIdMustExist)
# and row 40 (SARIF2001), whereas every other section break in the sheet
# uses a single blank row. Remove two of those blank rows so the SARIF2001+
# rows shift up to line up with the rest of the table (SARIF2001 -> row 38,
# SARIF2002 -> row 40, SARIF2003 -> row 42, etc.)
$ws.Rows("37:38").Delete()

# Provide messages for SARIF2003.ProvideVersionControlProvenance: its
# message-authoring/message-code status moves from TODO/TODO to
# READY FOR REVIEW/DONE (matching the green/blue fill used elsewhere).
$ws.Range("G2").Copy()
$ws.Range("G42").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("H42").PasteSpecial(-4122)
$ws.Range("G42").Value = "READY FOR REVIEW"
$ws.Range("H42").Value = "DONE"

# Match the saved view state (zoom level and active selection).
$excel.ActiveWindow.Zoom = 85
$ws.Range("H42").Select() | Out-Null
